$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an existing "Completed" cell as the formatting template so the
# cell style (font/fill matching the built-in "Good" cell style) is reused
# instead of Excel fabricating a brand-new style entry.
$template = $ws.Range("C4")
$template.Copy()

$cellsToComplete = @("H7", "E8", "H8", "C9", "E9")
foreach ($addr in $cellsToComplete) {
    $rng = $ws.Range($addr)
    $rng.PasteSpecial(-4122)  # xlPasteFormats
    $rng.Value = "Completed"
}

$excel.CutCopyMode = 0

# Update the active selection/cursor position on the sheet
$ws.Range("I5").Select()
